# Sprint 1 burndown chart: reduce remaining-hours figures for the last five
# days (22-27 Mar, columns L:Q) on Task 6 through Task 10 (rows 8-12).
# Row 13 ("Total Remaining") holds a SUM formula over each column and will
# recalculate automatically; the chart series (which read straight from the
# sheet) pick up the new cached values on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(14, 12, 10, 5, 2, 0)
$cols = @("L", "M", "N", "O", "P", "Q")

foreach ($row in 8..12) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $newValues[$i]
    }
}

# Leave the selection where the author left it.
$ws.Range("Q10").Select()
